# Apply updated cryptocurrency market data (prices / 1h volume changes)
# and a rank swap between TRON (row 12) and Polkadot (row 13).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    # Force the literal text into the cell even when it looks like a
    # number (e.g. "242.95" or "0.100") so trailing zeros / exact
    # formatting survive, then drop the temporary Text number format
    # so the cell keeps its original (default) style.
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

$ws.Range('D2').Value = '42.216.81'
$ws.Range('E2').Value = '  -0.33%  '
$ws.Range('D3').Value = '2.234.64'
$ws.Range('E3').Value = '  -0.02%  '
$ws.Range('E4').Value = '  +0.03%  '
Set-TextValue 'D5' '242.95'
Set-TextValue 'D6' '0.626'
$ws.Range('E6').Value = '  -0.26%  '
Set-TextValue 'D7' '74.08'
$ws.Range('E7').Value = '  -0.41%  '
$ws.Range('E8').Value = '  +0.17%  '
Set-TextValue 'D9' '0.600'
$ws.Range('E9').Value = '  -2.46%  '
Set-TextValue 'D10' '42.26'
$ws.Range('E10').Value = '  -2.23%  '
$ws.Range('E11').Value = '  +1.14%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D12' '6.93'
$ws.Range('E12').Value = '  -2.57%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue 'D13' '0.103'
$ws.Range('E13').Value = '  +0.37%  '
$ws.Range('D14').Value = '2.571.21'
Set-TextValue 'D15' '14.34'
$ws.Range('E15').Value = '  -0.56%  '
Set-TextValue 'D16' '0.837'
$ws.Range('E16').Value = '  -1.60%  '
$ws.Range('D17').Value = '2.214.69'
$ws.Range('E17').Value = '  -1.21%  '
$ws.Range('D18').Value = '42.082.59'
$ws.Range('E18').Value = '  -0.29%  '
$ws.Range('E19').Value = '  +0.73%  '
Set-TextValue 'D20' '6.22'
$ws.Range('E20').Value = '  +1.26%  '
Set-TextValue 'D21' '72.79'
$ws.Range('E21').Value = '  +1.29%  '
Set-TextValue 'D22' '11.47'
$ws.Range('E22').Value = '  +12.65%  '
Set-TextValue 'D23' '230.01'
$ws.Range('E23').Value = '  -0.19%  '
Set-TextValue 'D24' '2.03'
$ws.Range('E24').Value = '  -6.48%  '
$ws.Range('E25').Value = '  -0.11%  '
Set-TextValue 'D26' '11.37'
$ws.Range('E26').Value = '  -1.62%  '
$ws.Range('E27').Value = '  -0.72%  '
$ws.Range('E28').Value = '  -1.12%  '
$ws.Range('E29').Value = '  -2.09%  '
Set-TextValue 'D30' '167.21'
$ws.Range('E30').Value = '  +0.50%  '
Set-TextValue 'D31' '20.59'
$ws.Range('E31').Value = '  -1.49%  '
$ws.Range('E32').Value = '  -4.21%  '
$ws.Range('E33').Value = '  -0.36%  '
Set-TextValue 'D34' '29.83'
$ws.Range('E34').Value = '  +1.47%  '
$ws.Range('E35').Value = '  -0.10%  '
$ws.Range('E36').Value = '  -6.64%  '
Set-TextValue 'D37' '4.37'
$ws.Range('E37').Value = '  -3.63%  '
$ws.Range('E38').Value = '  -1.70%  '
Set-TextValue 'D39' '13.21'
$ws.Range('E39').Value = '  +0.22%  '
$ws.Range('E40').Value = '  -1.26%  '
$ws.Range('E41').Value = '  +0.71%  '
Set-TextValue 'D42' '64.54'
$ws.Range('E42').Value = '  +2.07%  '
$ws.Range('E43').Value = '  +0.29%  '
$ws.Range('E44').Value = '  -0.93%  '
Set-TextValue 'D45' '104.47'
$ws.Range('E45').Value = '  -1.14%  '
Set-TextValue 'D46' '0.100'
$ws.Range('E46').Value = '  -2.23%  '
$ws.Range('E47').Value = '  -0.12%  '
Set-TextValue 'D48' '2.34'
$ws.Range('E48').Value = '  -1.46%  '
$ws.Range('E49').Value = '  +0.16%  '
$ws.Range('E50').Value = '  -1.72%  '
$ws.Range('D51').Value = '2.443.91'
$ws.Range('E51').Value = '  +0.00%  '

Write-Output "Applied cryptos list update (82 cell changes)"
